$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add one point for interactivity: G12 changes from 0 to 1
$ws.Range("G12").Value = 1

# Update selection / view state to match the post-edit state
$ws.Range("G13").Select()
